# Applies the "1.0 -> 1.1" content update to the RF012 - Ajuda test-case sheet.
# The four repeating test-case blocks (TC1..TC4) keep their row positions,
# but the feature each block refers to is rotated:
#   TC1: Competencias (portfolio)  -> Periodos Avaliativos
#   TC2: Periodos Avaliativos      -> Avaliacoes
#   TC3: Avaliacoes                -> Niveis das Competencias
#   TC4: Niveis das Competencias   -> Competencias (portfolio)
#   TC5: Perfis de Competencias    -> unchanged

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC1 block (rows 10-11) now describes "Periodos Avaliativos"
$ws.Range("B10").Value = "Lider de Pessoas acessa a funcionalidade de gestao de Periodos Avaliativos a partir do menu inicial"
$ws.Range("D10").Value = "SYSTEM exibe a listagem dos Periodos Avaliativos cadastrados com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$ws.Range("B11").Value = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Periodos Avaliativos"

# TC2 block (rows 19-20) now describes "Avaliacoes"
$ws.Range("B19").Value = "Lider de Pessoas acessa a funcionalidade de gestao de Avaliacoes a partir do menu inicial"
$ws.Range("D19").Value = "SYSTEM exibe a listagem das Avaliacoes cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$ws.Range("B20").Value = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Avaliacoes"

# TC3 block (rows 28-29) now describes "Niveis das Competencias"
$ws.Range("B28").Value = "Lider de Pessoas acessa a funcionalidade de gestao de Niveis das Competencias a partir do menu inicial"
$ws.Range("D28").Value = "SYSTEM exibe a listagem dos Niveis das Competencias cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$ws.Range("B29").Value = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Niveis das Competencias"

# TC4 block (rows 37-38) now describes "Competencias (portfolio)"
$ws.Range("B37").Value = "Lider de Pessoas acessa a funcionalidade de gestao de Competencias (portfolio) a partir do menu inicial"
$ws.Range("D37").Value = "SYSTEM exibe a listagem das Competencias (portfolio) cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$ws.Range("B38").Value = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Competencias (portfolio)"

# TC5 block (rows 46-47, "Perfis de Competencias") is unchanged.
